$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 578
$ws.Range("I39").Value = 39.333332
$ws.Range("J39").Value = 1116.6666
$ws.Range("K39").Value = 117.999996
$ws.Range("L39").Value = 3349.9998
$ws.Range("M39").Value = 178.000004
$ws.Range("N39").Value = -3941.9998

$ws.Range("H40").Value = 2852.3333
$ws.Range("J40").Value = 3756.125
$ws.Range("L40").Value = 3756.125
$ws.Range("N40").Value = -4106.125

$ws.Range("H51").Value = 5166.6665
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 5500
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 5500
$ws.Range("M51").Value = -4516
$ws.Range("N51").Value = -6468

$ws.Range("H69").Value = 100000
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = $null

$ws.Range("H72").Value = 100000
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = $null

$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = $null
$ws.Range("N112").Value = $null

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").Value = $null

$ws.Range("H118").Value = 825
$ws.Range("I118").Value = 800
$ws.Range("J118").Value = 850
$ws.Range("K118").Value = 2400
$ws.Range("L118").Value = 2550
$ws.Range("M118").Value = -743
$ws.Range("N118").Value = -5864

$ws.Range("H125").Value = 630.1111
$ws.Range("I125").Value = 444.2
$ws.Range("K125").Value = 3997.8
$ws.Range("M125").Value = -1537.8

$ws.Range("H127").Value = 4000
$ws.Range("I127").Value = 4000
$ws.Range("K127").Value = 12000
$ws.Range("M127").Value = -7040

$ws.Range("H137").Value = 5045.1816
$ws.Range("I137").Value = 4500
$ws.Range("J137").Value = 5699.4
$ws.Range("K137").Value = 13500
$ws.Range("L137").Value = 17098.2
$ws.Range("M137").Value = -10950
$ws.Range("N137").Value = -22198.2

$ws.Range("H138").Value = 2494.5334
$ws.Range("J138").Value = 2594.25
$ws.Range("L138").Value = 7782.75
$ws.Range("N138").Value = -18062.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 998.5
$ws.Range("I110").Value = 998.5
$ws.Range("K110").Value = 998.5
$ws.Range("M110").Value = 1046.5

$ws.Range("H122").Value = 3942.077
$ws.Range("I122").Value = 5009.4
$ws.Range("J122").Value = 3275
$ws.Range("K122").Value = 15028.2
$ws.Range("L122").Value = 9825
$ws.Range("M122").Value = -12578.2
$ws.Range("N122").Value = -14725

$ws.Range("H132").Value = 3312.375
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 3285.5715
$ws.Range("K132").Value = 10500
$ws.Range("L132").Value = 9856.7145
$ws.Range("M132").Value = -7970
$ws.Range("N132").Value = -14916.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 9959.6
$ws.Range("I94").Value = 9999.5
$ws.Range("J94").Value = 9933
$ws.Range("K94").Value = 9999.5
$ws.Range("L94").Value = 9933
$ws.Range("M94").Value = -9548.5
$ws.Range("N94").Value = -10835

$ws.Range("H134").Value = 6865.0835
$ws.Range("I134").Value = 4480.5
$ws.Range("K134").Value = 13441.5
$ws.Range("M134").Value = -10906.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 735.9231
$ws.Range("I16").Value = 695.8182
$ws.Range("K16").Value = 695.8182
$ws.Range("M16").Value = -408.8182

$ws.Range("H22").Value = 2054.8
$ws.Range("I22").Value = 2546.5
$ws.Range("J22").Value = 88
$ws.Range("K22").Value = 2546.5
$ws.Range("L22").Value = 88
$ws.Range("M22").Value = -2196.5
$ws.Range("N22").Value = -788

$ws.Range("H62").Value = 2499
$ws.Range("I62").Value = 2499
$ws.Range("K62").Value = 2499
$ws.Range("M62").Value = -1875

$ws.Range("H65").Value = 2499
$ws.Range("I65").Value = 2499
$ws.Range("K65").Value = 12495
$ws.Range("M65").Value = -9375

$ws.Range("H68").Value = 35000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = $null

$ws.Range("H71").Value = 35000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = $null

$ws.Range("H113").Value = 735.9231
$ws.Range("I113").Value = 695.8182
$ws.Range("K113").Value = 695.8182
$ws.Range("M113").Value = 1474.1818

$ws.Range("H134").Value = 6000
$ws.Range("I134").Value = 6000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 18000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -15465
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1170.4286
$ws.Range("I68").Value = 1333
$ws.Range("K68").Value = 3999
$ws.Range("M68").Value = -3188

$ws.Range("H71").Value = 1170.4286
$ws.Range("I71").Value = 1333
$ws.Range("K71").Value = 11997
$ws.Range("M71").Value = -7941

$ws.Range("H92").Value = 1999.5
$ws.Range("I92").Value = 1999
$ws.Range("J92").Value = 2000
$ws.Range("K92").Value = 5997
$ws.Range("L92").Value = 6000
$ws.Range("M92").Value = -4749
$ws.Range("N92").Value = -8496

$ws.Range("H107").Value = 240
$ws.Range("J107").Value = 240
$ws.Range("L107").Value = 720
$ws.Range("N107").Value = -4560

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 9999
$ws.Range("J102").Value = 9999
$ws.Range("L102").Value = 9999
$ws.Range("N102").Value = -13243

$ws.Range("H122").Value = 715.8333
$ws.Range("J122").Value = 799
$ws.Range("L122").Value = 2397
$ws.Range("N122").Value = -7297

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1112.8572
$ws.Range("I16").Value = 1112.8572
$ws.Range("K16").Value = 1112.8572
$ws.Range("M16").Value = -942.8571999999999

$ws.Range("H132").Value = 4114.5713
$ws.Range("I132").Value = 4114.5713
$ws.Range("K132").Value = 12343.7139
$ws.Range("M132").Value = -9813.713899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1854860.9
$ws.Range("I2").Value = 2266941
$ws.Range("K2").Value = 2266941
$ws.Range("M2").Value = -2266829

$ws.Range("H4").Value = 31948.846
$ws.Range("I4").Value = 40689.2
$ws.Range("J4").Value = 2814.3333
$ws.Range("K4").Value = 40689.2
$ws.Range("L4").Value = 2814.3333
$ws.Range("M4").Value = -40576.2
$ws.Range("N4").Value = -3040.3333

$ws.Range("H54").Value = 30000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").Value = $null

$ws.Range("H81").Value = 37750.75
$ws.Range("I81").Value = 50001
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 100002
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -98941
$ws.Range("N81").Value = -4122

$ws.Range("H84").Value = 37750.75
$ws.Range("I84").Value = 50001
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 500010
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -494706
$ws.Range("N84").Value = -20608

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = $null

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = $null
